$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typo "IMMONOLOGISCHE" -> "IMMUNOLOGISCHE" in the ParameterCodeDisease
# column (F) for the disease rows 2-8.
$ws.Range("F2").Value = "COV_GECCO_IMMUNOLOGISCHE_ERKRANKUNG_CHRONISCH_ENTZUENDLICHE_DARMERKRANKUNG"
$ws.Range("F3").Value = "COV_GECCO_IMMUNOLOGISCHE_ERKRANKUNG_RHEUMATOIDE_ARTHRITIS"
$ws.Range("F4").Value = "COV_GECCO_IMMUNOLOGISCHE_ERKRANKUNG_KOLLAGENOSEN"
$ws.Range("F5").Value = "COV_GECCO_IMMUNOLOGISCHE_ERKRANKUNG_VASKULITIS"
$ws.Range("F6").Value = "COV_GECCO_IMMUNOLOGISCHE_ERKRANKUNG_ANGEBORENE_IMMUNDEFEKTE"
$ws.Range("F7").Value = "COV_GECCO_IMMUNOLOGISCHE_ERKRANKUNG_RHEUMATOLOGISCHE_ERKRANKUNG"
$ws.Range("F8").Value = "COV_GECCO_IMMUNOLOGISCHE_ERKRANKUNG_IMMUNOLOGISCHE_ERKRANKUNG"

# Move the active cell selection to F9, as saved by the author.
$ws.Range("F9").Select()
